$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 398.42856
$ws.Range("I4").Value = 399.5
$ws.Range("K4").Value = 399.5
$ws.Range("M4").Value = -285.5
$ws.Range("H99").Value = 587
$ws.Range("I99").Value = 587
$ws.Range("K99").Value = 1761
$ws.Range("M99").Value = -263
$ws.Range("H112").Value = 2882.9343
$ws.Range("J112").Value = 2882.9343
$ws.Range("L112").Value = 8648.802899999999
$ws.Range("N112").Value = -10864.8029
$ws.Range("H132").Value = 2315.1155
$ws.Range("I132").Value = 2270.5715
$ws.Range("J132").Value = 2502.2
$ws.Range("K132").Value = 6811.7145
$ws.Range("L132").Value = 7506.599999999999
$ws.Range("M132").Value = -4281.7145
$ws.Range("N132").Value = -12566.6
$ws.Range("H134").Value = 80780
$ws.Range("J134").Value = 80780
$ws.Range("L134").Value = 80780
$ws.Range("N134").Value = -90920
$ws.Range("H136").Value = 72500
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 72500
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 72500
$ws.Range("N136").Value = -82700
$ws.Range("H137").Value = 3802.8215
$ws.Range("J137").Value = 3339.1
$ws.Range("L137").Value = 10017.3
$ws.Range("N137").Value = -15117.3
$ws.Range("H138").Value = 4724.7144
$ws.Range("I138").Value = 945.28
$ws.Range("J138").Value = 10282.706
$ws.Range("K138").Value = 2835.84
$ws.Range("L138").Value = 30848.118
$ws.Range("M138").Value = 2304.16
$ws.Range("N138").Value = -41128.118
$ws.Range("M136").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 172.5
$ws.Range("I4").Value = 172.5
$ws.Range("K4").Value = 172.5
$ws.Range("M4").Value = -56.5
$ws.Range("H45").Value = 6880.9165
$ws.Range("I45").Value = 3821.375
$ws.Range("K45").Value = 3821.375
$ws.Range("M45").Value = -3444.375
$ws.Range("H61").Value = 6686.107
$ws.Range("I61").Value = 1900.0555
$ws.Range("K61").Value = 1900.0555
$ws.Range("M61").Value = -1688.0555
$ws.Range("H97").Value = 3087155
$ws.Range("I97").Value = 543.1429000000001
$ws.Range("J97").Value = 13890297
$ws.Range("K97").Value = 543.1429000000001
$ws.Range("L97").Value = 13890297
$ws.Range("M97").Value = -47.14290000000005
$ws.Range("N97").Value = -13891289
$ws.Range("H132").Value = 4056.9343
$ws.Range("I132").Value = 1901.5143
$ws.Range("K132").Value = 5704.5429
$ws.Range("M132").Value = -3174.5429
$ws.Range("H136").Value = 6686.107
$ws.Range("I136").Value = 1900.0555
$ws.Range("K136").Value = 5700.166499999999
$ws.Range("M136").Value = -3150.166499999999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2792.5217
$ws.Range("I105").Value = 2190.625
$ws.Range("J105").Value = 4168.2856
$ws.Range("K105").Value = 2190.625
$ws.Range("L105").Value = 4168.2856
$ws.Range("M105").Value = -443.625
$ws.Range("N105").Value = -7662.2856
$ws.Range("H107").Value = 45003644
$ws.Range("I107").Value = 53573770
$ws.Range("K107").Value = 53573770
$ws.Range("M107").Value = -53571850
$ws.Range("H128").Value = 3547.7778
$ws.Range("I128").Value = 3547.7778
$ws.Range("K128").Value = 10643.3334
$ws.Range("M128").Value = -8153.3334
$ws.Range("H134").Value = 6151.275
$ws.Range("I134").Value = 2102.4348
$ws.Range("J134").Value = 11629.117
$ws.Range("K134").Value = 6307.3044
$ws.Range("L134").Value = 34887.351
$ws.Range("M134").Value = -3772.3044
$ws.Range("N134").Value = -39957.351

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 12201962
$ws.Range("I58").Value = 35717224
$ws.Range("J58").Value = 8864.406999999999
$ws.Range("K58").Value = 35717224
$ws.Range("L58").Value = 8864.406999999999
$ws.Range("M58").Value = -35717021
$ws.Range("N58").Value = -9270.406999999999
$ws.Range("H62").Value = 16164.833
$ws.Range("I62").Value = 3663.3333
$ws.Range("K62").Value = 3663.3333
$ws.Range("M62").Value = -3039.3333
$ws.Range("H65").Value = 16164.833
$ws.Range("I65").Value = 3663.3333
$ws.Range("K65").Value = 18316.6665
$ws.Range("M65").Value = -15196.6665
$ws.Range("H76").Value = 4992.1665
$ws.Range("I76").Value = 4992.1665
$ws.Range("K76").Value = 4992.1665
$ws.Range("M76").Value = -4677.1665
$ws.Range("H79").Value = 4992.1665
$ws.Range("I79").Value = 4992.1665
$ws.Range("K79").Value = 4992.1665
$ws.Range("M79").Value = -3900.1665
$ws.Range("H132").Value = 8476.24
$ws.Range("I132").Value = 4683.375
$ws.Range("K132").Value = 14050.125
$ws.Range("M132").Value = -11520.125
$ws.Range("H134").Value = 5487.3335
$ws.Range("I134").Value = 2629.652
$ws.Range("J134").Value = 7607.5483
$ws.Range("K134").Value = 7888.956
$ws.Range("L134").Value = 22822.6449
$ws.Range("M134").Value = -5353.956
$ws.Range("N134").Value = -27892.6449
$ws.Range("H136").Value = 12201962
$ws.Range("I136").Value = 35717224
$ws.Range("J136").Value = 8864.406999999999
$ws.Range("K136").Value = 107151672
$ws.Range("L136").Value = 26593.221
$ws.Range("M136").Value = -107149122
$ws.Range("N136").Value = -31693.221

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 16667959
$ws.Range("I14").Value = 16667959
$ws.Range("K14").Value = 50003877
$ws.Range("M14").Value = -50003704
$ws.Range("H87").Value = 250002260
$ws.Range("I87").Value = 333333980
$ws.Range("K87").Value = 1000001940
$ws.Range("M87").Value = -1000000692
$ws.Range("H90").Value = 250002260
$ws.Range("I90").Value = 333333980
$ws.Range("K90").Value = 3000005820
$ws.Range("M90").Value = -2999999580
$ws.Range("H113").Value = 3145.5217
$ws.Range("J113").Value = 4069.3572
$ws.Range("L113").Value = 12208.0716
$ws.Range("N113").Value = -16548.0716
$ws.Range("H140").Value = 1820.6364
$ws.Range("I140").Value = 973.8823
$ws.Range("J140").Value = 4699.6
$ws.Range("K140").Value = 2921.6469
$ws.Range("L140").Value = 14098.8
$ws.Range("M140").Value = 2258.3531
$ws.Range("N140").Value = -24458.8

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 93634.45
$ws.Range("I80").Value = 2996
$ws.Range("J80").Value = 202400.6
$ws.Range("K80").Value = 2996
$ws.Range("L80").Value = 202400.6
$ws.Range("M80").Value = -1998
$ws.Range("N80").Value = -204396.6
$ws.Range("H83").Value = 93634.45
$ws.Range("I83").Value = 2996
$ws.Range("J83").Value = 202400.6
$ws.Range("K83").Value = 14980
$ws.Range("L83").Value = 1012003
$ws.Range("M83").Value = -9988
$ws.Range("N83").Value = -1021987
$ws.Range("H132").Value = 10232.228
$ws.Range("I132").Value = 5365.4287
$ws.Range("K132").Value = 16096.2861
$ws.Range("M132").Value = -13566.2861

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2384.0557
$ws.Range("J22").Value = 4573.2856
$ws.Range("L22").Value = 4573.2856
$ws.Range("N22").Value = -5163.2856
$ws.Range("H27").Value = 2384.0557
$ws.Range("J27").Value = 4573.2856
$ws.Range("L27").Value = 4573.2856
$ws.Range("N27").Value = -4787.2856
$ws.Range("H82").Value = 743397.6
$ws.Range("I82").Value = 1761990.4
$ws.Range("J82").Value = 2602.9092
$ws.Range("K82").Value = 1761990.4
$ws.Range("L82").Value = 2602.9092
$ws.Range("M82").Value = -1761629.4
$ws.Range("N82").Value = -3324.9092
$ws.Range("H85").Value = 743397.6
$ws.Range("I85").Value = 1761990.4
$ws.Range("J85").Value = 2602.9092
$ws.Range("K85").Value = 1761990.4
$ws.Range("L85").Value = 2602.9092
$ws.Range("M85").Value = -1760742.4
$ws.Range("N85").Value = -5098.9092
$ws.Range("H107").Value = 2543.7778
$ws.Range("I107").Value = 2543.7778
$ws.Range("K107").Value = 2543.7778
$ws.Range("M107").Value = -623.7777999999998
$ws.Range("H132").Value = 10210987
$ws.Range("I132").Value = 21742258
$ws.Range("K132").Value = 65226774
$ws.Range("M132").Value = -65224244
$ws.Range("H136").Value = 5941.9297
$ws.Range("I136").Value = 1685.8611
$ws.Range("J136").Value = 13238.048
$ws.Range("K136").Value = 5057.5833
$ws.Range("L136").Value = 39714.144
$ws.Range("M136").Value = -2507.5833
$ws.Range("N136").Value = -44814.144

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5132.5557
$ws.Range("I62").Value = 6038.6
$ws.Range("K62").Value = 6038.6
$ws.Range("M62").Value = -5414.6
$ws.Range("H65").Value = 5132.5557
$ws.Range("I65").Value = 6038.6
$ws.Range("K65").Value = 30193
$ws.Range("M65").Value = -27073
$ws.Range("H107").Value = 13334089
$ws.Range("I107").Value = 422.76923
$ws.Range("J107").Value = 27778894
$ws.Range("K107").Value = 1268.30769
$ws.Range("L107").Value = 83336682
$ws.Range("M107").Value = 651.6923099999999
$ws.Range("N107").Value = -83340522
$ws.Range("H132").Value = 62516576
$ws.Range("I132").Value = 100021720
$ws.Range("J132").Value = 8000
$ws.Range("K132").Value = 300065160
$ws.Range("L132").Value = 24000
$ws.Range("M132").Value = -300062630
$ws.Range("N132").Value = -29060
$ws.Range("H136").Value = 19253442
$ws.Range("I136").Value = 33334238
$ws.Range("K136").Value = 100002714
$ws.Range("M136").Value = -100000164
